$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 gained a new budget line item: "Parts for Prototyping" at $151.16
$ws.Range("A15").Value = "Parts for Prototyping"
$ws.Range("B15").Value = 151.16

# The active selection in the saved file moved to B16
[void]$ws.Range("B16").Select()
